# Powerpoint writer: consolidate text runs when possible.
# Merge adjacent a:r runs that share identical formatting into a single
# run by rewriting the covering TextRange's .Text (leaves differently
# formatted runs, like the Consolas "inline code"/"head" spans, alone).

$p = $ppt.ActivePresentation

# Slide 1 title: "Header" + " " + "with" + " " + "inline code"(Consolas)
#   -> "Header with " + "inline code"(Consolas)
$tr1 = $p.Slides.Item(1).Shapes.Item(1).TextFrame.TextRange
$sub1 = $tr1.Characters(1, 12)
$sub1.Text = "Header with "

# Slide 2 title: "Syntax" + " " + "highlighting" -> "Syntax highlighting"
$tr2 = $p.Slides.Item(2).Shapes.Item(1).TextFrame.TextRange
$sub2 = $tr2.Characters(1, $tr2.Length)
$sub2.Text = "Syntax highlighting"

# Slide 3 title: "Two" + " " + "column" + " " + "slide" -> "Two column slide"
$tr3 = $p.Slides.Item(3).Shapes.Item(1).TextFrame.TextRange
$sub3 = $tr3.Characters(1, $tr3.Length)
$sub3.Text = "Two column slide"
